# "add comparison to whisper and glossy, add latency images captured on sat"
#
# 1) Fix the Ron-time-Whisper value on Sheet1 (4.1 -> 4), which ripples
#    through the dependent formulas in row 31.
# 2) Add a new "Sheet2" after Sheet1 with a Whisper/Glossy comparison table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: correct the Whisper Ron-time figure -------------------------
$ws1.Range("C31").Value = 4

# --- Add Sheet2 right after Sheet1 ---------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate the text/header cells FIRST, and in this precise order, so the
# shared-string table is built up in the same order as the authored file.
$ws2.Range("B3").Value  = "ma"
$ws2.Range("A4").Value  = "Whisper (ms)"
$ws2.Range("A5").Value  = "Glossy (ms)"
$ws2.Range("E3").Value  = "hours"
$ws2.Range("E6").Value  = "hour"
$ws2.Range("C3").Value  = "active (ms)"
$ws2.Range("D3").Value  = "days"
$ws2.Range("F3").Value  = "battery capacity (mAh)"
$ws2.Range("G10").Value = "mAh"
$ws2.Range("F10").Value = "Ron time Whisper"
$ws2.Range("B16").Value = "Formula is: (Packet size * 8) x TTL + 500"
$ws2.Range("I11").Value = "(I suspect actual is 2.1"
$ws2.Range("I13").Value = "I suspect actual is 4.1"

# Remaining header cell without a 1:1 new string (D6 re-uses "hours")
$ws2.Range("D6").Value = "hours"

# --- Row 4: Whisper ---------------------------------------------------
$ws2.Range("B4").Value   = 20
$ws2.Range("C4").Value   = 1.9
$ws2.Range("D4").Value   = 2193
$ws2.Range("E4").Formula = "=D4*24"
$ws2.Range("F4").Value   = 2000
$ws2.Range("G4").Formula = "=2000/(B4*(2/3600000))"

# --- Row 5: Glossy ------------------------------------------------------
$ws2.Range("B5").Value   = 20
$ws2.Range("C5").Value   = 3.7
$ws2.Range("D5").Value   = 1126
$ws2.Range("E5").Formula = "=D5*24"
$ws2.Range("F5").Value   = 2000

# --- Row 7 ---------------------------------------------------------------
$ws2.Range("B7").Value = 20
$ws2.Range("C7").Value = 1.9
$ws2.Range("D7").Value = 52632
$ws2.Range("E7").Value = 3600

# --- Row 8 (scientific number format) -------------------------------------
$ws2.Range("B8").Formula = "=20 *(2/3600000)"
$ws2.Range("B8").NumberFormat = "0.00E+00"
$ws2.Range("C8").Formula = "=2000/B8"

# --- Row 9 -----------------------------------------------------------------
$ws2.Range("F9").Formula = "=2*360000"

# --- Rows 11-14: Ron-time derivation table ---------------------------------
$ws2.Range("E11").Formula = "=C4"
$ws2.Range("F11").Formula = "=((E11/1000)*360)"
$ws2.Range("G11").Formula = "=ROUND(F11/20,3)"
$ws2.Range("H11").Formula = "=2000/G11"

$ws2.Range("E12").Formula = "=C5"
$ws2.Range("F12").Value   = 0.76
$ws2.Range("G12").Formula = "=F12/20"
$ws2.Range("H12").Formula = "=2000/G12"

$ws2.Range("E13").Formula = "=C5"
$ws2.Range("F13").Formula = "=((E13/1000)*360)"
$ws2.Range("G13").Formula = "=ROUND(F13/20,3)"
$ws2.Range("H13").Formula = "=2000/G13"

$ws2.Range("E14").Value   = 1.2
$ws2.Range("F14").Formula = "=((E14/1000)*360)"
$ws2.Range("G14").Formula = "=F14/20"
$ws2.Range("H14").Formula = "=2000/G14"

# --- Column widths (best-effort match of the authored layout) -------------
$ws2.Columns.Item(2).ColumnWidth = 74.422
$ws2.Columns.Item(3).ColumnWidth = 11.422
$ws2.Columns.Item(5).ColumnWidth = 13.2541
$ws2.Columns.Item(6).ColumnWidth = 19.09
$ws2.Columns.Item(7).ColumnWidth = 11.2541

# --- Selections / active views --------------------------------------------
[void]$ws1.Range("C34").Select()
$w = $excel.ActiveWindow
$w.ScrollRow    = 4
$w.ScrollColumn = 1

[void]$ws2.Range("D29").Select()
[void]$ws2.Activate()
